$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B2: the JSON-like key changes from "requirement" to "request"
# (text content otherwise stays the same).
$ws.Range("B2").Value = '"request": "Lauther wants the average density for the months of January and February for the measurement system with tag F980-40 on the platform Reconvavo"'

# Update the view: scroll so row 2 is at top, and select B4 as the active cell.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
